# BOT; UPDATE DATA
# Appends two new daily rows (2020-05-28 / 2020-05-29) to the "相談件数"
# sheet just above the trailing footnote row, pushing the footnote from
# row 124 down to row 126, and refreshes the sheet's print area / view
# state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# The footnote currently sits on row 124 ("A124" blank / "B124" the
# shared-string footnote). Insert two blank rows above it (copying the
# formatting of the row above, same as Excel's native Insert), which
# shoves the footnote down to row 126 and frees up rows 124-125 for the
# new data.
$ws.Rows.Item(124).Insert()
$ws.Rows.Item(124).Insert()

# New data row for 2020-05-28 (date serial 43979).
$ws.Range("A124").Value = 43979
$ws.Range("B124").Value = 124
$ws.Range("C124").Value = 39431
$ws.Range("D124").Value = 27
$ws.Range("E124").Value = 7939

# New data row for 2020-05-29 (date serial 43980).
$ws.Range("A125").Value = 43980
$ws.Range("B125").Value = 139
$ws.Range("C125").Value = 39570
$ws.Range("D125").Value = 24
$ws.Range("E125").Value = 7963

# Extend the print area to cover the two new rows.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$126"
    }
}

# Update the saved selection/scroll state to land on the new last cell.
$ws.Range("F125").Select()
